# Daily attendance processing - 2025-10-31 21:41:26
# Reorders the "Recorded By" (column G) values so that the literal token
# "System" (exact case) — when present among the comma-separated list of
# recorders — is moved to the front of the list, leaving the relative
# order of all other tokens (including a differently-cased "system")
# unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$col = 7  # column G = "Recorded By"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -notmatch ",") { continue }

    $parts = $val -split ",\s*"

    $idx = -1
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($parts[$i].Equals("System")) {
            $idx = $i
            break
        }
    }

    if ($idx -le 0) { continue }

    $newVal = "System"
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($i -ne $idx) {
            $newVal = $newVal + ", " + $parts[$i]
        }
    }

    $cell.Value = $newVal
}
